$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the dates in column F (rows 2-7) forward by one day
$ws.Range("F2").Value = 44596
$ws.Range("F3").Value = 44595
$ws.Range("F4").Value = 44594
$ws.Range("F5").Value = 44593
$ws.Range("F6").Value = 44592
$ws.Range("F7").Value = 44591
